# Append a new data row (row 4) to the materials sheet:
#   A4 = "" (blank, matching the existing blank cells in column A)
#   B4 = "CSC103"
#   C4 = "store/materials/CSC103/conditional probability.pdf"
#
# The worksheet's used range grows from A1:C3 to A1:C4 as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Touching the style (without actually changing it) materializes the cell
# in the sheet even though its value is blank - mirroring the existing
# A2/A3 placeholder cells in column A.
$ws.Range("A4").Style = "Normal"

$ws.Range("B4").Value = "CSC103"
$ws.Range("C4").Value = "store/materials/CSC103/conditional probability.pdf"
